# Update cryptos list data (price and volume columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.278.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.506.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.504.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("E11").Value = "  -4.46%  "

$ws.Range("E12").Value = "  -2.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.102.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("E14").Value = "  -3.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.510.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.309.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("E23").Value = "  -1.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.640.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -1.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.54%  "

$ws.Range("E29").Value = "  -3.14%  "

$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("E33").Value = "  -8.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.494.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("E37").Value = "  -3.49%  "

$ws.Range("E38").Value = "  -4.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.68%  "

$ws.Range("E44").Value = "  -1.52%  "

$ws.Range("E45").Value = "  -9.11%  "

$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.942"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.24%  "
